$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.424.87"
$ws.Range("E2").Value = "  +1.25%  "

Set-TextValue $ws.Range("D3") "2.684.57"
$ws.Range("E3").Value = "  +2.74%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.16%  "

Set-TextValue $ws.Range("D5") "534.35"
$ws.Range("E5").Value = "  +4.13%  "

Set-TextValue $ws.Range("D6") "157.56"
$ws.Range("E6").Value = "  +2.39%  "

$ws.Range("E7").Value = "  -0.21%  "

Set-TextValue $ws.Range("D8") "0.592"
$ws.Range("E8").Value = "  +1.05%  "

Set-TextValue $ws.Range("D9") "6.62"
$ws.Range("E9").Value = "  -2.75%  "

Set-TextValue $ws.Range("D10") "0.110"
$ws.Range("E10").Value = "  +5.11%  "

Set-TextValue $ws.Range("D11") "0.355"
$ws.Range("E11").Value = "  +2.82%  "

$ws.Range("E12").Value = "  -0.07%  "

Set-TextValue $ws.Range("D13") "3.148.88"
$ws.Range("E13").Value = "  +2.54%  "

Set-TextValue $ws.Range("D14") "61.389.81"
$ws.Range("E14").Value = "  +1.33%  "

Set-TextValue $ws.Range("D15") "22.28"
$ws.Range("E15").Value = "  +3.07%  "

$ws.Range("E16").Value = "  +2.84%  "

Set-TextValue $ws.Range("D17") "2.693.56"
$ws.Range("E17").Value = "  +2.79%  "

Set-TextValue $ws.Range("D18") "4.81"
$ws.Range("E18").Value = "  +1.51%  "

Set-TextValue $ws.Range("D19") "357.80"
$ws.Range("E19").Value = "  +1.13%  "

Set-TextValue $ws.Range("D20") "10.80"
$ws.Range("E20").Value = "  +2.15%  "

Set-TextValue $ws.Range("D21") "6.39"
$ws.Range("E21").Value = "  +3.59%  "

Set-TextValue $ws.Range("D22") "0.999"
$ws.Range("E22").Value = "  -0.02%  "

Set-TextValue $ws.Range("D23") "62.10"
$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("E24").Value = "  +3.20%  "

Set-TextValue $ws.Range("D25") "0.170"
$ws.Range("E25").Value = "  +2.22%  "

Set-TextValue $ws.Range("D26") "0.997"
$ws.Range("E26").Value = "  +0.18%  "

Set-TextValue $ws.Range("D27") "0.0₃0875"
$ws.Range("E27").Value = "  +3.96%  "

Set-TextValue $ws.Range("D28") "7.46"
$ws.Range("E28").Value = "  +1.67%  "

Set-TextValue $ws.Range("D29") "0.998"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  +6.35%  "

Set-TextValue $ws.Range("D31") "19.70"
$ws.Range("E31").Value = "  +1.50%  "

Set-TextValue $ws.Range("D32") "1.64"
$ws.Range("E32").Value = "  +4.02%  "

Set-TextValue $ws.Range("D33") "150.30"
$ws.Range("E33").Value = "  -1.05%  "

Set-TextValue $ws.Range("D34") "4.19"
$ws.Range("E34").Value = "  +5.37%  "

Set-TextValue $ws.Range("D35") "1.22"
$ws.Range("E35").Value = "  +2.76%  "

Set-TextValue $ws.Range("D36") "0.924"
$ws.Range("E36").Value = "  +9.91%  "

Set-TextValue $ws.Range("D37") "0.891"
$ws.Range("E37").Value = "  +2.95%  "

Set-TextValue $ws.Range("D38") "1.52"
$ws.Range("E38").Value = "  +2.24%  "

$ws.Range("E41").Value = "  +1.84%  "

Set-TextValue $ws.Range("D42") "0.659"
$ws.Range("E42").Value = "  +5.45%  "

Set-TextValue $ws.Range("D43") "21.09"
$ws.Range("E43").Value = "  +6.60%  "

Set-TextValue $ws.Range("D44") "0.103"
$ws.Range("E44").Value = "  +1.38%  "

Set-TextValue $ws.Range("D45") "0.0571"
$ws.Range("E45").Value = "  +3.23%  "

$ws.Range("E46").Value = "  +0.04%  "

Set-TextValue $ws.Range("D47") "5.06"
$ws.Range("E47").Value = "  +2.87%  "

Set-TextValue $ws.Range("D48") "0.0243"
$ws.Range("E48").Value = "  +3.87%  "

Set-TextValue $ws.Range("D49") "19.27"
$ws.Range("E49").Value = "  +9.33%  "

# Row 39 <-> Row 40 swap (Bittensor/OKB)
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D39") "36.95"
$ws.Range("E39").Value = "  +1.54%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D40") "307.50"
$ws.Range("E40").Value = "  +5.18%  "

# Row 50 <-> Row 51 swap (WhiteBITCoin/Maker)
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D50") "2.020.23"
$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D51") "10.34"
$ws.Range("E51").Value = "  +0.25%  "
